# Auto-generated edit script: updates Leve profit-calculation cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 729.6667
$ws.Range("I29").Value = 729.6667
$ws.Range("K29").Value = 2189.0001
$ws.Range("M29").Value = -1908.0001
$ws.Range("H33").Value = 339.94446
$ws.Range("I33").Value = 322.85715
$ws.Range("K33").Value = 322.85715
$ws.Range("M33").Value = -93.85714999999999
$ws.Range("H99").Value = 641.6
$ws.Range("I99").Value = 641.6
$ws.Range("K99").Value = 1924.8
$ws.Range("M99").Value = -426.8000000000002
$ws.Range("H129").Value = 883.42
$ws.Range("I129").Value = 430.66666
$ws.Range("J129").Value = 912.31915
$ws.Range("K129").Value = 1291.99998
$ws.Range("L129").Value = 2736.95745
$ws.Range("M129").Value = 3708.00002
$ws.Range("N129").Value = -12736.95745
$ws.Range("H137").Value = 3765.5
$ws.Range("I137").Value = 2367.5833
$ws.Range("K137").Value = 7102.749899999999
$ws.Range("M137").Value = -4552.749899999999
$ws.Range("H138").Value = 3551.697
$ws.Range("I138").Value = 1185.2778
$ws.Range("J138").Value = 4077.5679
$ws.Range("K138").Value = 3555.8334
$ws.Range("L138").Value = 12232.7037
$ws.Range("M138").Value = 1584.1666
$ws.Range("N138").Value = -22512.7037

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6308.691
$ws.Range("I32").Value = 5432.282
$ws.Range("J32").Value = 8444.9375
$ws.Range("K32").Value = 5432.282
$ws.Range("L32").Value = 8444.9375
$ws.Range("M32").Value = -5145.282
$ws.Range("N32").Value = -9018.9375
$ws.Range("H61").Value = 1733.4445
$ws.Range("I61").Value = 1250.1666
$ws.Range("K61").Value = 1250.1666
$ws.Range("M61").Value = -1038.1666
$ws.Range("H74").Value = 7910.769
$ws.Range("I74").Value = 9201.333000000001
$ws.Range("J74").Value = 5007
$ws.Range("K74").Value = 9201.333000000001
$ws.Range("L74").Value = 5007
$ws.Range("M74").Value = -8327.333000000001
$ws.Range("N74").Value = -6755
$ws.Range("H77").Value = 7910.769
$ws.Range("I77").Value = 9201.333000000001
$ws.Range("J77").Value = 5007
$ws.Range("K77").Value = 46006.665
$ws.Range("L77").Value = 25035
$ws.Range("M77").Value = -41638.665
$ws.Range("N77").Value = -33771
$ws.Range("H110").Value = 1180.1904
$ws.Range("I110").Value = 1215.3158
$ws.Range("K110").Value = 1215.3158
$ws.Range("M110").Value = 829.6841999999999
$ws.Range("H122").Value = 2108.5908
$ws.Range("I122").Value = 1595.8235
$ws.Range("J122").Value = 3852
$ws.Range("K122").Value = 4787.470499999999
$ws.Range("L122").Value = 11556
$ws.Range("M122").Value = -2337.470499999999
$ws.Range("N122").Value = -16456
$ws.Range("H136").Value = 1733.4445
$ws.Range("I136").Value = 1250.1666
$ws.Range("K136").Value = 3750.4998
$ws.Range("M136").Value = -1200.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3194.889
$ws.Range("J99").Value = 5119.9
$ws.Range("L99").Value = 5119.9
$ws.Range("N99").Value = -8115.9
$ws.Range("H134").Value = 2024.3
$ws.Range("I134").Value = 1513.5143
$ws.Range("J134").Value = 5599.8
$ws.Range("K134").Value = 4540.5429
$ws.Range("L134").Value = 16799.4
$ws.Range("M134").Value = -2005.5429
$ws.Range("N134").Value = -21869.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14708642
$ws.Range("I31").Value = 1370.7368
$ws.Range("K31").Value = 1370.7368
$ws.Range("M31").Value = -1075.7368
$ws.Range("H34").Value = 14708642
$ws.Range("I34").Value = 1370.7368
$ws.Range("K34").Value = 1370.7368
$ws.Range("M34").Value = -1168.7368
$ws.Range("H105").Value = 3666.6667
$ws.Range("I105").Value = 3000
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 3000
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -1253
$ws.Range("N105").Value = -8494
$ws.Range("H132").Value = 1515.8182
$ws.Range("I132").Value = 1044.1052
$ws.Range("J132").Value = 4503.3335
$ws.Range("K132").Value = 3132.3156
$ws.Range("L132").Value = 13510.0005
$ws.Range("M132").Value = -602.3155999999999
$ws.Range("N132").Value = -18570.0005
$ws.Range("H134").Value = 4983.8667
$ws.Range("I134").Value = 5599.6
$ws.Range("J134").Value = 3752.4
$ws.Range("K134").Value = 16798.8
$ws.Range("L134").Value = 11257.2
$ws.Range("M134").Value = -14263.8
$ws.Range("N134").Value = -16327.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 38910
$ws.Range("I4").Value = 150150
$ws.Range("K4").Value = 450450
$ws.Range("M4").Value = -450338
$ws.Range("H45").Value = 2266.6667
$ws.Range("J45").Value = 2850
$ws.Range("L45").Value = 8550
$ws.Range("N45").Value = -9614
$ws.Range("H69").Value = 10833.333
$ws.Range("I69").Value = 1666.6666
$ws.Range("J69").Value = 20000
$ws.Range("K69").Value = 4999.9998
$ws.Range("L69").Value = 60000
$ws.Range("M69").Value = -4188.9998
$ws.Range("N69").Value = -61622
$ws.Range("H72").Value = 10833.333
$ws.Range("I72").Value = 1666.6666
$ws.Range("J72").Value = 20000
$ws.Range("K72").Value = 14999.9994
$ws.Range("L72").Value = 180000
$ws.Range("M72").Value = -10943.9994
$ws.Range("N72").Value = -188112
$ws.Range("H131").Value = 10204970
$ws.Range("I131").Value = 100000410
$ws.Range("J131").Value = 942.0682
$ws.Range("K131").Value = 300001230
$ws.Range("L131").Value = 2826.2046
$ws.Range("M131").Value = -299996190
$ws.Range("N131").Value = -12906.2046

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1179.75
$ws.Range("I97").Value = 1107.6
$ws.Range("K97").Value = 1107.6
$ws.Range("M97").Value = -611.5999999999999
$ws.Range("H113").Value = 2066.4119
$ws.Range("I113").Value = 2552
$ws.Range("J113").Value = 1520.125
$ws.Range("K113").Value = 2552
$ws.Range("L113").Value = 1520.125
$ws.Range("M113").Value = -382
$ws.Range("N113").Value = -5860.125
$ws.Range("H132").Value = 3525.5
$ws.Range("I132").Value = 1757.0625
$ws.Range("J132").Value = 7062.375
$ws.Range("K132").Value = 5271.1875
$ws.Range("L132").Value = 21187.125
$ws.Range("M132").Value = -2741.1875
$ws.Range("N132").Value = -26247.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 46695
$ws.Range("J63").Value = 46695
$ws.Range("L63").Value = 46695
$ws.Range("N63").Value = -48193
$ws.Range("H66").Value = 46695
$ws.Range("J66").Value = 46695
$ws.Range("L66").Value = 140085
$ws.Range("N66").Value = -147573
$ws.Range("H93").Value = 4277365.5
$ws.Range("I93").Value = 7411074
$ws.Range("J93").Value = 4126.364
$ws.Range("K93").Value = 7411074
$ws.Range("L93").Value = 4126.364
$ws.Range("M93").Value = -7409826
$ws.Range("N93").Value = -6622.364
$ws.Range("H122").Value = 2624.4167
$ws.Range("I122").Value = 1693.3334
$ws.Range("J122").Value = 5417.6665
$ws.Range("K122").Value = 5080.0002
$ws.Range("L122").Value = 16252.9995
$ws.Range("M122").Value = -2630.0002
$ws.Range("N122").Value = -21152.9995
$ws.Range("H136").Value = 3481.1785
$ws.Range("I136").Value = 1511.4117
$ws.Range("J136").Value = 6525.364
$ws.Range("K136").Value = 4534.2351
$ws.Range("L136").Value = 19576.092
$ws.Range("M136").Value = -1984.2351
$ws.Range("N136").Value = -24676.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1881.1428
$ws.Range("I136").Value = 1225.2
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 3675.6
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -1125.6
$ws.Range("N136").Value = -50100

